# Apply the commit "Add files via upload" changes to Sheet1:
#  - C2 (板单保存99%无响应 row) keeps the spsavebaneditlog writeup, but the
#    thread count in bullet 2 goes from 100 to 1000.
#  - C4 (PDA出库异常 row) gets a more detailed sp_pdaOutCaseCommit writeup
#    (adds a C# exception-handling bullet under "建议：").
#  - Row 4's height grows (27 -> 67.5) to fit the longer text, and the
#    sheet's last active-cell selection moves from C4 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$banText = "1. 修改存储过程spsavebaneditlog捕获异常，发现ban_makebill_itemlist_log表字段remotematname长度为50，而ban_makebill_itemlist相应字段长度为500，造成“截断字符串”的错误，在未捕获异常的情况下，前端是否有处理？建议：修改spsavebaneditlog捕获异常并抛出错误信息`n2. 多线程压力测试，1000条线程同时调用spsavebaneditlog，在不发生“截断字符串”的错误的情况下，并未发现任何异常。`n3. 归档日志表"

$pdaText = "1. 分析存储过程sp_pdaOutCaseCommit`n建议： `n   - SET SACT_ABORT ON`n   - C#应用程序加错误捕捉和处理语句（处理超时异常）`n2. 分析“General Network Error,Check your Network Documentation”"

$ws.Range("C2").Value = $banText
$ws.Range("C4").Value = $pdaText

$ws.Rows.Item(4).RowHeight = 67.5

$ws.Range("C3").Select()
